$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 544; everything from row 544 down shifts to 545+.
$ws.Rows("544:544").Insert()

# Populate the newly inserted row 544 with the new weekly price record.
$ws.Range("A544").Value = 10
$ws.Range("B544").Value = "Vega Modelo de Temuco"
$ws.Range("C544").Value = "La Araucanía"
$ws.Range("D544").Value = 45223
$ws.Range("E544").Value = 9
$ws.Range("F544").Value = 100112017
$ws.Range("G544").Value = "Apio"
$ws.Range("H544").Value = "Americana (o)"
$ws.Range("I544").Value = "Primera"
$ws.Range("J544").Value = 400
$ws.Range("K544").Value = 8000
$ws.Range("L544").Value = 8000
$ws.Range("M544").Value = 8000
$ws.Range("N544").Value = "`$/caja 8 unidades"
$ws.Range("O544").Value = "Provincia del Elquí"
$ws.Range("P544").Value = 8000
$ws.Range("Q544").Value = 1
$ws.Range("R544").Value = "Hortaliza"
